$wb = $excel.ActiveWorkbook

# Translate the worksheet tab name from English to Vietnamese
$ws = $wb.Worksheets.Item("ERoute")
$ws.Name = "Tuyến"
